# Updated cryptos list on Thu Nov  2 21:49:10 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold plain numeric-looking text (e.g. "233.05", "35.119.72").
# Force Text format before assigning so Excel does not auto-convert them to
# Number cells, then restore the default "Normal" style so no stray number
# format is left behind on the cell (matches original unstyled inline-string cells).
function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "35.119.72"
$ws.Range("E2").Value = "  -0.39%  "
Set-TextValue "D3" "1.813.46"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  +0.67%  "
Set-TextValue "D5" "233.05"
$ws.Range("E5").Value = "  +2.17%  "
Set-TextValue "D6" "0.612"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.70%  "
$ws.Range("E8").Value = "  -6.09%  "
Set-TextValue "D9" "0.323"
$ws.Range("E9").Value = "  +5.90%  "
$ws.Range("E10").Value = "  -1.37%  "
$ws.Range("E11").Value = "  -0.48%  "
Set-TextValue "D12" "2.076.90"
$ws.Range("E12").Value = "  -1.61%  "
Set-TextValue "D13" "1.810.80"
$ws.Range("E13").Value = "  -1.64%  "
$ws.Range("E14").Value = "  +0.49%  "
Set-TextValue "D15" "11.05"
$ws.Range("E15").Value = "  -4.55%  "
$ws.Range("E16").Value = "  -1.52%  "
Set-TextValue "D17" "35.080.65"
$ws.Range("E17").Value = "  -0.37%  "
Set-TextValue "D18" "69.56"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("E19").Value = "  -0.47%  "
Set-TextValue "D20" "238.66"
$ws.Range("E20").Value = "  -2.90%  "
Set-TextValue "D21" "11.89"
$ws.Range("E21").Value = "  -1.93%  "
Set-TextValue "D22" "4.70"
$ws.Range("E22").Value = "  -0.27%  "
$ws.Range("E23").Value = "  +0.72%  "
$ws.Range("E24").Value = "  +3.22%  "
Set-TextValue "D25" "172.04"
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("E26").Value = "  -1.11%  "
Set-TextValue "D27" "17.50"
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("E28").Value = "  -1.43%  "
$ws.Range("E29").Value = "  +21.50%  "
$ws.Range("E30").Value = "  +0.74%  "
Set-TextValue "D31" "4.20"
$ws.Range("E31").Value = "  +6.56%  "
Set-TextValue "D32" "3.329.05"
$ws.Range("E32").Value = "  -8.70%  "
Set-TextValue "D33" "0.0552"
$ws.Range("E33").Value = "  +2.70%  "
$ws.Range("E34").Value = "  -0.86%  "
$ws.Range("E35").Value = "  -5.87%  "
$ws.Range("E36").Value = "  +5.07%  "
Set-TextValue "D37" "92.41"
$ws.Range("E37").Value = "  +2.24%  "
$ws.Range("E38").Value = "  +0.53%  "
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D40" "1.310.18"
$ws.Range("E40").Value = "  -2.39%  "
$ws.Range("B41").Value = "WEMIXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D41" "1.28"
$ws.Range("E41").Value = "  +1.49%  "
$ws.Range("E42").Value = "  -1.99%  "
Set-TextValue "D43" "2.48"
$ws.Range("E43").Value = "  +1.10%  "
Set-TextValue "D44" "14.54"
$ws.Range("E44").Value = "  -2.27%  "
$ws.Range("E45").Value = "  -5.79%  "
Set-TextValue "D46" "2.77"
$ws.Range("E46").Value = "  -2.12%  "
$ws.Range("E47").Value = "  +4.12%  "
Set-TextValue "D48" "0.0512"
$ws.Range("E48").Value = "  -1.56%  "
Set-TextValue "D49" "1.990.94"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("E50").Value = "  +0.68%  "
Set-TextValue "D51" "0.0651"
$ws.Range("E51").Value = "  +4.88%  "

Write-Host "Updated cryptos list cells"
